$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.07837893064909
$ws.Range("D2").Value = 1.079566417411833
$ws.Range("E2").Value = 1.081813542453623
$ws.Range("F2").Value = 1.09241558652532
$ws.Range("I2").Value = 1.065828316021834
$ws.Range("J2").Value = 1.083268700750428
$ws.Range("K2").Value = 1.082242132426751
$ws.Range("L2").Value = 1.084483386601136
$ws.Range("M2").Value = 1.095058083056891
$ws.Range("N2").Value = 1.08480706655815

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.07959198113016
$ws.Range("D3").Value = 1.080532148426952
$ws.Range("E3").Value = 1.082879779012723
$ws.Range("F3").Value = 1.093514371441239
$ws.Range("I3").Value = 1.066269136261808
$ws.Range("J3").Value = 1.08414108101649
$ws.Range("K3").Value = 1.083025598248837
$ws.Range("L3").Value = 1.085367538662868
$ws.Range("M3").Value = 1.095976683625121
$ws.Range("N3").Value = 1.085680685704253

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.08037678007734
$ws.Range("D4").Value = 1.08115686553006
$ws.Range("E4").Value = 1.083569740833363
$ws.Range("F4").Value = 1.094225499648296
$ws.Range("I4").Value = 1.066553056129498
$ws.Range("J4").Value = 1.084704862235279
$ws.Range("K4").Value = 1.083531752630428
$ws.Range("L4").Value = 1.085939073949416
$ws.Range("M4").Value = 1.096570618829723
$ws.Range("N4").Value = 1.08624526755704

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.080706680972378
$ws.Range("D5").Value = 1.081419455194067
$ws.Range("E5").Value = 1.083859810034468
$ws.Range("F5").Value = 1.094524492402239
$ws.Range("I5").Value = 1.06667210045177
$ws.Range("J5").Value = 1.084941707423531
$ws.Range("K5").Value = 1.08374434900166
$ws.Range("L5").Value = 1.086179211637558
$ws.Range("M5").Value = 1.096820199596786
$ws.Range("N5").Value = 1.086482449092626

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.080762071153088
$ws.Range("D6").Value = 1.081463542743335
$ws.Range("E6").Value = 1.083908514494366
$ws.Range("F6").Value = 1.094574696642889
$ws.Range("I6").Value = 1.066692070025327
$ws.Range("J6").Value = 1.084981464920572
$ws.Range("K6").Value = 1.083780033691539
$ws.Range("L6").Value = 1.086219523874116
$ws.Range("M6").Value = 1.096862098900188
$ws.Range("N6").Value = 1.086522263049872

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.080381188339614
$ws.Range("D7").Value = 1.081160374427953
$ws.Range("E7").Value = 1.083573616714915
$ws.Range("F7").Value = 1.094229494667806
$ws.Range("I7").Value = 1.066554648044991
$ws.Range("J7").Value = 1.084708027633039
$ws.Range("K7").Value = 1.083534594102912
$ws.Range("L7").Value = 1.085942283213386
$ws.Range("M7").Value = 1.096573954169867
$ws.Range("N7").Value = 1.086248437450028

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.078788913719164
$ws.Range("D8").Value = 1.079892827414135
$ws.Range("E8").Value = 1.082173875204106
$ws.Range("F8").Value = 1.092786897296179
$ws.Range("I8").Value = 1.065977566960274
$ws.Range("J8").Value = 1.083563672680722
$ws.Range("K8").Value = 1.082507074508706
$ws.Range("L8").Value = 1.084782308457934
$ws.Range("M8").Value = 1.095368624262695
$ws.Range("N8").Value = 1.085102457382418

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.075982089545532
$ws.Range("D9").Value = 1.077657876914791
$ws.Range("E9").Value = 1.079707589075979
$ws.Range("F9").Value = 1.090245902586715
$ws.Range("I9").Value = 1.064950542150429
$ws.Range("J9").Value = 1.081541718517676
$ws.Range("K9").Value = 1.080690292667388
$ws.Range("L9").Value = 1.082733882989944
$ws.Range("M9").Value = 1.093241117278717
$ws.Range("N9").Value = 1.08307763181254

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.074110080002922
$ws.Range("D10").Value = 1.076166947884789
$ws.Range("E10").Value = 1.078063507114903
$ws.Range("F10").Value = 1.088552560761606
$ws.Range("I10").Value = 1.064259012773203
$ws.Range("J10").Value = 1.080190025563587
$ws.Range("K10").Value = 1.079474917004553
$ws.Range("L10").Value = 1.081365255973167
$ws.Range("M10").Value = 1.091820336634268
$ws.Range("N10").Value = 1.081724019299415

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.073299268251385
$ws.Range("D11").Value = 1.075521121146257
$ws.Range("E11").Value = 1.077351615161485
$ws.Range("F11").Value = 1.087819468618517
$ws.Range("I11").Value = 1.063957941409895
$ws.Range("J11").Value = 1.07960383121581
$ws.Range("K11").Value = 1.078947641847343
$ws.Range("L11").Value = 1.080771899455348
$ws.Range("M11").Value = 1.091204532125644
$ws.Range("N11").Value = 1.081136992488426

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.072998061663421
$ws.Range("D12").Value = 1.075281194639636
$ws.Range("E12").Value = 1.077087186607276
$ws.Range("F12").Value = 1.087547184881869
$ws.Range("I12").Value = 1.063845863789657
$ws.Range("J12").Value = 1.079385955729217
$ws.Range("K12").Value = 1.078751635692242
$ws.Range("L12").Value = 1.080551389326021
$ws.Range("M12").Value = 1.090975704119054
$ws.Range("N12").Value = 1.080918807593652

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.073062673116007
$ws.Range("D13").Value = 1.07533266141102
$ws.Range("E13").Value = 1.077143907456581
$ws.Range("F13").Value = 1.087605589810998
$ws.Range("I13").Value = 1.063869915977504
$ws.Range("J13").Value = 1.079432696971235
$ws.Range("K13").Value = 1.078793686603604
$ws.Range("L13").Value = 1.080598694542875
$ws.Range("M13").Value = 1.091024792629248
$ws.Range("N13").Value = 1.080965615213595

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.073274371146806
$ws.Range("D14").Value = 1.075501289520939
$ws.Range("E14").Value = 1.07732975742512
$ws.Range("F14").Value = 1.087796961172096
$ws.Range("I14").Value = 1.063948682067747
$ws.Range("J14").Value = 1.079585824358047
$ws.Range("K14").Value = 1.078931443038056
$ws.Range("L14").Value = 1.080753674292179
$ws.Range("M14").Value = 1.091185618989665
$ws.Range("N14").Value = 1.08111896005886

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.073404800590773
$ws.Range("D15").Value = 1.075605181840835
$ws.Range("E15").Value = 1.077444265665684
$ws.Range("F15").Value = 1.087914873911333
$ws.Range("I15").Value = 1.063997179795068
$ws.Range("J15").Value = 1.079680153029947
$ws.Range("K15").Value = 1.079016299061049
$ws.Range("L15").Value = 1.080849147675255
$ws.Range("M15").Value = 1.091284697354822
$ws.Range("N15").Value = 1.081213422688293

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.074163886113033
$ws.Range("D16").Value = 1.076209804043532
$ws.Range("E16").Value = 1.078110753076109
$ws.Range("F16").Value = 1.088601216444079
$ws.Range("I16").Value = 1.064278959386682
$ws.Range("J16").Value = 1.08022891023452
$ws.Range("K16").Value = 1.079509889164447
$ws.Range("L16").Value = 1.081404619572788
$ws.Range("M16").Value = 1.091861192878409
$ws.Range("N16").Value = 1.081762959191041

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.074639980249187
$ws.Range("D17").Value = 1.076589001415
$ws.Range("E17").Value = 1.078528823922223
$ws.Range("F17").Value = 1.08903177661968
$ws.Range("I17").Value = 1.06445527413874
$ws.Range("J17").Value = 1.08057288866846
$ws.Range("K17").Value = 1.079819234209354
$ws.Range("L17").Value = 1.081752855650004
$ws.Range("M17").Value = 1.092222652660766
$ws.Range("N17").Value = 1.082107426113812

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.074917656886131
$ws.Range("D18").Value = 1.076810157286717
$ws.Range("E18").Value = 1.078772678266516
$ws.Range("F18").Value = 1.089282928286378
$ws.Range("I18").Value = 1.064557957886546
$ws.Range("J18").Value = 1.080773438485616
$ws.Range("K18").Value = 1.079999572609279
$ws.Range("L18").Value = 1.081955905311519
$ws.Range("M18").Value = 1.09243342849633
$ws.Range("N18").Value = 1.082308260734712

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.075012333941179
$ws.Range("D19").Value = 1.076885561759865
$ws.Range("E19").Value = 1.078855826448499
$ws.Range("F19").Value = 1.089368566787268
$ws.Range("I19").Value = 1.064592943670727
$ws.Range("J19").Value = 1.080841806074152
$ws.Range("K19").Value = 1.08006104686535
$ws.Range("L19").Value = 1.082025128066467
$ws.Range("M19").Value = 1.092505287882476
$ws.Range("N19").Value = 1.082376725413065

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.074588902042668
$ws.Range("D20").Value = 1.076548319567197
$ws.Range("E20").Value = 1.078483968813414
$ws.Range("F20").Value = 1.088985580230035
$ws.Range("I20").Value = 1.064436373547586
$ws.Range("J20").Value = 1.080535992043726
$ws.Range("K20").Value = 1.07978605448827
$ws.Range("L20").Value = 1.081715500519928
$ws.Range("M20").Value = 1.092183877406221
$ws.Range("N20").Value = 1.082070477091639

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.073212032335942
$ws.Range("D21").Value = 1.07545163379385
$ws.Range("E21").Value = 1.077275029255901
$ws.Range("F21").Value = 1.087740606551502
$ws.Range("I21").Value = 1.063925494222632
$ws.Range("J21").Value = 1.079540735945639
$ws.Range("K21").Value = 1.07889088142466
$ws.Range("L21").Value = 1.080708039694713
$ws.Range("M21").Value = 1.091138262150062
$ws.Range("N21").Value = 1.081073807615734

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0723461361725
$ws.Range("D22").Value = 1.074761885326947
$ws.Range("E22").Value = 1.076514919331263
$ws.Range("F22").Value = 1.086957953450677
$ws.Range("I22").Value = 1.063602858673951
$ws.Range("J22").Value = 1.078914186809707
$ws.Range("K22").Value = 1.078327166696601
$ws.Range("L22").Value = 1.080073965743134
$ws.Range("M22").Value = 1.090480316498514
$ws.Range("N22").Value = 1.080446368708164

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.072805184092519
$ws.Range("D23").Value = 1.07512755507693
$ws.Range("E23").Value = 1.076917868522494
$ws.Range("F23").Value = 1.087372842493212
$ws.Range("I23").Value = 1.063774029242886
$ws.Range("J23").Value = 1.079246407857799
$ws.Range("K23").Value = 1.078626086664396
$ws.Range("L23").Value = 1.080410161699945
$ws.Range("M23").Value = 1.090829156142902
$ws.Range("N23").Value = 1.08077906154825

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.07461198214995
$ws.Range("D24").Value = 1.076566702015505
$ws.Range("E24").Value = 1.078504236903817
$ws.Range("F24").Value = 1.089006454348869
$ws.Range("I24").Value = 1.064444914398089
$ws.Range("J24").Value = 1.080552664308642
$ws.Range("K24").Value = 1.079801047276439
$ws.Range("L24").Value = 1.081732379912919
$ws.Range("M24").Value = 1.092201398452135
$ws.Range("N24").Value = 1.082087173033083

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.076707853467083
$ws.Range("D25").Value = 1.078235831784271
$ws.Range("E25").Value = 1.080345160260596
$ws.Range("F25").Value = 1.090902691663162
$ws.Range("I25").Value = 1.065217256986436
$ws.Range("J25").Value = 1.08206509430955
$ws.Range("K25").Value = 1.081160709112574
$ws.Range("L25").Value = 1.083263976908189
$ws.Range("M25").Value = 1.093791555846008
$ws.Range("N25").Value = 1.083601750858071
